$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "temp_in"
$ws.Range("D1").Value = "temp_out"
$ws.Range("E1").Value = "press_in"
$ws.Range("F1").Value = "press_out"
$ws.Range("G1").Value = "volume"
$ws.Range("H1").Value = "steps"
$ws.Range("I1").Value = "molar_flow_in"

$ws.Range("I2").Select()
